# Fill in the "Соответствие" (correspondence) column for the
# "Критическая технология" (KT) and "Сквозная технология" (ST) tables.
#
# Word's Tables collection is 1-indexed; in $d.Tables these are table
# index 2 (KT) and table index 3 (ST) -- table index 1 is the earlier
# "Приоритетная проблема медицины и здравоохранения" table.

$d = $word.ActiveDocument

$ktTable = $d.Tables.Item(2)
$ktTable.Cell(2, 2).Range.Text = "Нет"
$ktTable.Cell(3, 2).Range.Text = "Да"
$ktTable.Cell(4, 2).Range.Text = "Нет"
$ktTable.Cell(5, 2).Range.Text = "Нет"

$stTable = $d.Tables.Item(3)
$stTable.Cell(2, 2).Range.Text = "Да"
$stTable.Cell(3, 2).Range.Text = "Нет"
$stTable.Cell(4, 2).Range.Text = "Нет"
$stTable.Cell(5, 2).Range.Text = "Нет"
